# Handback status report regeneration:
#   - File "4a1a6148-7426-402f-b4cd-6613294a232b.md" was handed back again at a later
#     timestamp and is now identified as "76c58e28-84cf-41ee-b692-2bb5c64269f6.md"
#   - File "ee03db45-3948-4fe5-a1cb-902e4795ba7e.md" was handed back again at a later
#     timestamp and is now identified as "ffffea95f99f-e464-439d-805b-ad2129f8bb03.md"
#   - Both files ended up sharing the same xliff hash/content in zh-cn and de-de

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------------
# Sheet "Overview"
# ----------------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "76c58e28-84cf-41ee-b692-2bb5c64269f6.md"
$ov.Range("B2").Value = "e2e\76c58e28-84cf-41ee-b692-2bb5c64269f6.md"
$ov.Range("G2").Value = "2016-08-29 01:01:09"

$ov.Range("A3").Value = "ffffea95f99f-e464-439d-805b-ad2129f8bb03.md"
$ov.Range("B3").Value = "e2e\ffffea95f99f-e464-439d-805b-ad2129f8bb03.md"
$ov.Range("G3").Value = "2016-08-29 01:01:09"

# Rebuild the hyperlinks (same target addresses, updated display text)
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6ab939b926161d86e4d338e012f20a3db3d60509/e2e/4a1a6148-7426-402f-b4cd-6613294a232b.md", "", "", "e2e\76c58e28-84cf-41ee-b692-2bb5c64269f6.md")
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6ab939b926161d86e4d338e012f20a3db3d60509/e2e/ee03db45-3948-4fe5-a1cb-902e4795ba7e.md", "", "", "e2e\ffffea95f99f-e464-439d-805b-ad2129f8bb03.md")

# ----------------------------------------------------------------------------------
# Sheet "zh-cn"
# ----------------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "76c58e28-84cf-41ee-b692-2bb5c64269f6.md"
$zh.Range("G2").Value = "76c58e28-84cf-41ee-b692-2bb5c64269f6.a650dce89da627fa1327405f9da5a0e0a35480f4.zh-cn.xlf"
$zh.Range("H2").Value = "2016-08-29 01:01:01"
$zh.Range("I2").Value = "76c58e28-84cf-41ee-b692-2bb5c64269f6.md"
$zh.Range("J2").Value = "76c58e28-84cf-41ee-b692-2bb5c64269f6.a650dce89da627fa1327405f9da5a0e0a35480f4.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-29 01:01:42"

$zh.Range("A3").Value = "ffffea95f99f-e464-439d-805b-ad2129f8bb03.md"
$zh.Range("G3").Value = "76c58e28-84cf-41ee-b692-2bb5c64269f6.a650dce89da627fa1327405f9da5a0e0a35480f4.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-29 01:01:01"
$zh.Range("I3").Value = "ffffea95f99f-e464-439d-805b-ad2129f8bb03.md"
$zh.Range("J3").Value = "76c58e28-84cf-41ee-b692-2bb5c64269f6.a650dce89da627fa1327405f9da5a0e0a35480f4.zh-cn.xlf"
$zh.Range("K3").Value = "2016-08-29 01:01:42"

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6ab939b926161d86e4d338e012f20a3db3d60509/e2e/4a1a6148-7426-402f-b4cd-6613294a232b.md", "", "", "76c58e28-84cf-41ee-b692-2bb5c64269f6.md")
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/81f1a82db8d0740beaa69ecd796a78a3311a5ef2/e2e/4a1a6148-7426-402f-b4cd-6613294a232b.md", "", "", "76c58e28-84cf-41ee-b692-2bb5c64269f6.md")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6ab939b926161d86e4d338e012f20a3db3d60509/e2e/ee03db45-3948-4fe5-a1cb-902e4795ba7e.md", "", "", "ffffea95f99f-e464-439d-805b-ad2129f8bb03.md")
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/81f1a82db8d0740beaa69ecd796a78a3311a5ef2/e2e/ee03db45-3948-4fe5-a1cb-902e4795ba7e.md", "", "", "ffffea95f99f-e464-439d-805b-ad2129f8bb03.md")

# ----------------------------------------------------------------------------------
# Sheet "de-de"
# ----------------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "76c58e28-84cf-41ee-b692-2bb5c64269f6.md"
$de.Range("G2").Value = "76c58e28-84cf-41ee-b692-2bb5c64269f6.a650dce89da627fa1327405f9da5a0e0a35480f4.de-de.xlf"
$de.Range("H2").Value = "2016-08-29 01:01:09"
$de.Range("I2").Value = "76c58e28-84cf-41ee-b692-2bb5c64269f6.md"
$de.Range("J2").Value = "76c58e28-84cf-41ee-b692-2bb5c64269f6.a650dce89da627fa1327405f9da5a0e0a35480f4.de-de.xlf"
$de.Range("K2").Value = "2016-08-29 01:01:49"

$de.Range("A3").Value = "ffffea95f99f-e464-439d-805b-ad2129f8bb03.md"
$de.Range("G3").Value = "76c58e28-84cf-41ee-b692-2bb5c64269f6.a650dce89da627fa1327405f9da5a0e0a35480f4.de-de.xlf"
$de.Range("H3").Value = "2016-08-29 01:01:09"
$de.Range("I3").Value = "ffffea95f99f-e464-439d-805b-ad2129f8bb03.md"
$de.Range("J3").Value = "76c58e28-84cf-41ee-b692-2bb5c64269f6.a650dce89da627fa1327405f9da5a0e0a35480f4.de-de.xlf"
$de.Range("K3").Value = "2016-08-29 01:01:49"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6ab939b926161d86e4d338e012f20a3db3d60509/e2e/4a1a6148-7426-402f-b4cd-6613294a232b.md", "", "", "76c58e28-84cf-41ee-b692-2bb5c64269f6.md")
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/931045667e04b8aa616ce33aeea806a5c4780052/e2e/4a1a6148-7426-402f-b4cd-6613294a232b.md", "", "", "76c58e28-84cf-41ee-b692-2bb5c64269f6.md")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6ab939b926161d86e4d338e012f20a3db3d60509/e2e/ee03db45-3948-4fe5-a1cb-902e4795ba7e.md", "", "", "ffffea95f99f-e464-439d-805b-ad2129f8bb03.md")
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/931045667e04b8aa616ce33aeea806a5c4780052/e2e/ee03db45-3948-4fe5-a1cb-902e4795ba7e.md", "", "", "ffffea95f99f-e464-439d-805b-ad2129f8bb03.md")

Write-Host "Handback status report regenerated."
